# Apply weekly fruit/vegetable data corrections (Higo - Vega Central Mapocho de Santiago)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44320
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range("S2").Value = 1714

# Row 3
$ws.Range("D3").Value = 44320
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("S3").Value = 1143

# Row 8
$ws.Range("D8").Value = 44322
$ws.Range("M8").Value = 45
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 1714

# Row 9
$ws.Range("D9").Value = 44322
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 1143

# Row 10
$ws.Range("D10").Value = 44980
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 16000
$ws.Range("S10").Value = 2286

# Row 11
$ws.Range("D11").Value = 44980
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 13000
$ws.Range("O11").Value = 13000
$ws.Range("P11").Value = 13000
$ws.Range("S11").Value = 1857

# Row 12
$ws.Range("D12").Value = 44302
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("S12").Value = 2143

# Row 13
$ws.Range("D13").Value = 44302
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("S13").Value = 1714

# Row 15
$ws.Range("D15").Value = 44299
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 15000
$ws.Range("R15").Value = "Provincia de Santiago"
$ws.Range("S15").Value = 2143

# Row 16
$ws.Range("D16").Value = 44299
$ws.Range("M16").Value = 75
$ws.Range("N16").Value = 12000
$ws.Range("O16").Value = 12000
$ws.Range("P16").Value = 12000
$ws.Range("R16").Value = "Provincia de Santiago"
$ws.Range("S16").Value = 1714

# Row 17
$ws.Range("D17").Value = 44292
$ws.Range("M17").Value = 25
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("S17").Value = 2286

# Row 18
$ws.Range("D18").Value = 44292
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 15000
$ws.Range("S18").Value = 2143
